$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue "D2" '29.885.49'
Set-TextValue "E2" '  +0.93%  '
Set-TextValue "D3" '1.623.67'
Set-TextValue "E3" '  +1.06%  '
Set-TextValue "E4" '  -0.47%  '
Set-TextValue "D5" '214.29'
Set-TextValue "E5" '  +0.79%  '
Set-TextValue "D6" '0.522'
Set-TextValue "E6" '  +0.02%  '
Set-TextValue "E7" '  -0.43%  '
Set-TextValue "D8" '30.03'
Set-TextValue "E8" '  +11.75%  '
Set-TextValue "D9" '0.259'
Set-TextValue "E9" '  +2.82%  '
Set-TextValue "D10" '0.0611'
Set-TextValue "E10" '  +1.53%  '
Set-TextValue "E11" '  +0.44%  '
Set-TextValue "D12" '1.855.33'
Set-TextValue "E12" '  +0.94%  '
Set-TextValue "D13" '1.619.30'
Set-TextValue "E13" '  +1.28%  '
Set-TextValue "D14" '0.569'
Set-TextValue "E14" '  +5.99%  '
Set-TextValue "E15" '  +4.84%  '
Set-TextValue "D16" '29.929.93'
Set-TextValue "E16" '  +1.01%  '
Set-TextValue "D17" '8.83'
Set-TextValue "E17" '  +16.36%  '
Set-TextValue "D18" '64.62'
Set-TextValue "E18" '  +1.88%  '
Set-TextValue "D19" '244.05'
Set-TextValue "E19" '  +1.41%  '
Set-TextValue "E20" '  +1.52%  '
Set-TextValue "D21" '0.996'
Set-TextValue "E21" '  -0.30%  '
Set-TextValue "D22" '4.12'
Set-TextValue "E22" '  +3.17%  '
Set-TextValue "D23" '9.62'
Set-TextValue "E23" '  +4.17%  '
Set-TextValue "E24" '  +2.29%  '
Set-TextValue "D25" '157.11'
Set-TextValue "E25" '  +1.66%  '
Set-TextValue "E26" '  +2.38%  '
Set-TextValue "E27" '  +1.90%  '
Set-TextValue "E28" '  +2.95%  '
Set-TextValue "D29" '0.996'
Set-TextValue "E29" '  -0.43%  '
Set-TextValue "E30" '  +2.84%  '
Set-TextValue "E31" '  +4.98%  '
Set-TextValue "E32" '  +3.54%  '
Set-TextValue "E33" '  +3.21%  '
Set-TextValue "D34" '1.427.18'
Set-TextValue "E34" '  +1.05%  '
Set-TextValue "E35" '  +7.03%  '
Set-TextValue "E36" '  -0.22%  '
Set-TextValue "E37" '  +1.87%  '
Set-TextValue "E38" '  -0.65%  '
Set-TextValue "E39" '  +2.77%  '
Set-TextValue "D40" '0.558'
Set-TextValue "D41" '0.0505'
Set-TextValue "E41" '  +2.81%  '
Set-TextValue "E42" '  +0.61%  '
Set-TextValue "D43" '0.833'
Set-TextValue "E43" '  +4.31%  '
Set-TextValue "D44" '54.01'
Set-TextValue "E44" '  -0.04%  '
Set-TextValue "D45" '69.16'
Set-TextValue "E45" '  +4.74%  '
Set-TextValue "E46" '  +17.41%  '
Set-TextValue "E47" '  -0.39%  '
Set-TextValue "D48" '5.41'
Set-TextValue "E48" '  +2.31%  '
Set-TextValue "D49" '1.763.97'
Set-TextValue "E49" '  +0.89%  '
Set-TextValue "D50" '88.48'
Set-TextValue "E50" '  +2.17%  '
Set-TextValue "D51" '0.0₆0107'
Set-TextValue "E51" '  +3.21%  '
